# Update "view study program" sheet:
#  - Shorten the "knowledge type" labels in column O (drop the parenthetical
#    credit-count suffixes) and let Excel re-derive the shared-string table.
#  - Because the labels got shorter, the previously-manual row heights
#    (27.6pt, needed to wrap the long text) are no longer required, so we
#    auto-fit those rows back to the sheet's default height.
#  - Move the frozen-pane view / selection down toward the bottom of the
#    "Kiến thức cơ sở ngành" block (row 46) like the author was doing when
#    they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Column O label rewrites (credit counts dropped from the text) ----
$ws.Range("O2:O25").Value   = "Đại cương"
$ws.Range("O26:O31").Value  = "Đại cương không tích luỹ"
$ws.Range("O32:O47").Value  = "Kiến thức cơ sở ngành"
$ws.Range("O48:O53").Value  = "Chuyên ngành Công nghệ Phần mềm"
$ws.Range("O54:O59").Value  = "Chuyên ngành Công nghệ Dữ liệu"
$ws.Range("O60:O65").Value  = "Chuyên ngành An ninh Mạng và IoT"
$ws.Range("O66:O83").Value  = "Các môn tự chọn chuyên ngành"

# ---- Row heights: rows 26-83 no longer need the tall, wrapped height ----
$ws.Rows("26:83").AutoFit()

# ---- View state: scroll / select near the bottom of the edited block ----
$win = $excel.ActiveWindow
$ws.Range("A1").Select()
$win.FreezePanes = $false
$ws.Range("A2").Select()
$win.FreezePanes = $true
$ws.Range("O46").Select()
